$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") values from 45203 to 45205 for rows 2 through 295.
$ws.Range("C2:C295").Value = 45205
